$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("draftpicks")
$ws.Activate()

# New draft pick rows (445-461) — "more picks, refreshed projections"
$picks = @(
    @{ row = 445; team = "ds9";       player = "Albert Pujols";         salary = 1; position = "1B"; drafted = 43875 },
    @{ row = 446; team = "marmaduke"; player = "Wade Davis";             salary = 7; position = "P";  drafted = 43875 },
    @{ row = 447; team = "dembums";   player = "Austin Romine";          salary = 5; position = "C";  drafted = 43875 },
    @{ row = 448; team = "ottawa";    player = "Jake Bauers";            salary = 1; position = "CI"; drafted = 43875 },
    @{ row = 449; team = "drjames";   player = "Dexter Fowler";          salary = 2; position = "DH"; drafted = 43875 },
    @{ row = 450; team = "rippe";     player = "Mike Fiers";             salary = 1; position = "P";  drafted = 43875 },
    @{ row = 451; team = "isotopes";  player = "Daniel Hudson";          salary = 1; position = "P";  drafted = 43875 },
    @{ row = 452; team = "deano";     player = "Rich Hill";              salary = 1; position = "P";  drafted = 43875 },
    @{ row = 453; team = "drjames";   player = "Chad Green";             salary = 2; position = "P";  drafted = 43876 },
    @{ row = 454; team = "drjames";   player = "Seranthony Dominguez";   salary = 1; position = "P";  drafted = 43876 },
    @{ row = 455; team = "bears";     player = "Andres Munoz";           salary = 1; position = "P";  drafted = 43876 },
    @{ row = 456; team = "marmaduke"; player = "Stephen Vogt";           salary = 1; position = "C";  drafted = 43876 },
    @{ row = 457; team = "pasadena";  player = "Luis Urias";             salary = 1; position = "P";  drafted = 43876 },
    @{ row = 458; team = "dembums";   player = "Jordan Montgomery";      salary = 1; position = "P";  drafted = 43876 },
    @{ row = 459; team = "ds9";       player = "Ryan Mountcastle";       salary = 3; position = "CI"; drafted = 43876 },
    @{ row = 460; team = "sturgeon";  player = "Rowan Wick";             salary = 1; position = "P";  drafted = 43876 },
    @{ row = 461; team = "marmaduke"; player = "Todd Frazier";           salary = 2; position = "1B"; drafted = 43877 }
)

foreach ($p in $picks) {
    $r = $p.row
    $ws.Range("A$r").Value = $p.team
    $ws.Range("B$r").Value = $p.player
    $ws.Range("C$r").Value = $p.salary
    $ws.Range("D$r").Value = $p.position
    $ws.Range("E$r").Value = $p.drafted
}

# Stray formatted (but empty) cell left at A463, carrying the salary's
# currency number format, matching the source workbook's final state.
$ws.Range("A463").NumberFormat = "\$#,##0;[Red]\$#,##0"

# View state: zoomed in, scrolled down near the new rows, active cell on
# the last new pick.
$excel.ActiveWindow.Zoom = 115
$ws.Range("A461").Select() | Out-Null
